# Refresh the cryptocurrency price ("D") and 1h-volume-change ("E") columns
# with the latest scraped figures.
#
# Column D values are stored as plain text, not numbers (prices are shown
# with locale-style "." thousands separators, and some are round numbers like
# "1.000" where a trailing zero must be preserved). Excel normally infers a
# numeric type for text input like "236.34", which would silently reformat/
# truncate it, so each D-column write temporarily forces the Text ("@") number
# format, assigns the literal string, then clears the temporary formatting so
# the cell is left exactly as it started (General format, default style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "30.224.65"
$c.ClearFormats()
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.860.70"
$c.ClearFormats()
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "236.34"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.34%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.ClearFormats()
$ws.Range("E6").Value = "  +0.00%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4681"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.26%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2875"
$c.ClearFormats()
$ws.Range("E8").Value = "  +1.17%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06550"
$c.ClearFormats()
$ws.Range("E9").Value = "  +0.27%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "21.75"
$c.ClearFormats()
$ws.Range("E10").Value = "  +2.68%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07930"
$c.ClearFormats()
$ws.Range("E11").Value = "  +0.26%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "97.77"
$c.ClearFormats()
$ws.Range("E12").Value = "  +0.49%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.871.05"
$c.ClearFormats()
$ws.Range("E13").Value = "  -0.19%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.173"
$c.ClearFormats()
$ws.Range("E14").Value = "  +0.36%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6818"
$c.ClearFormats()
$ws.Range("E15").Value = "  +0.60%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "268.53"
$c.ClearFormats()
$ws.Range("E16").Value = "  -5.79%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "30.221.27"
$c.ClearFormats()
$ws.Range("E17").Value = "  -0.40%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "13.77"
$c.ClearFormats()
$ws.Range("E18").Value = "  +8.43%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.0000"
$c.ClearFormats()
$ws.Range("E19").Value = "  -0.06%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.000007454"
$c.ClearFormats()
$ws.Range("E20").Value = "  +2.45%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "2.111.63"
$c.ClearFormats()
$ws.Range("E21").Value = "  -0.56%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.327"
$c.ClearFormats()
$ws.Range("E23").Value = "  +0.15%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.186"
$c.ClearFormats()
$ws.Range("E24").Value = "  -0.16%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "167.18"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("E26").Value = "  -0.90%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.93"
$c.ClearFormats()
$ws.Range("E27").Value = "  -1.02%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.963"
$c.ClearFormats()
$ws.Range("E28").Value = "  +2.05%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.386"
$c.ClearFormats()
$ws.Range("E29").Value = "  +2.48%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.09841"
$c.ClearFormats()
$ws.Range("E30").Value = "  +1.69%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.393"
$c.ClearFormats()
$ws.Range("E31").Value = "  -0.95%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.074"
$c.ClearFormats()
$ws.Range("E33").Value = "  -0.81%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.04714"
$c.ClearFormats()
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  +1.36%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7037"
$c.ClearFormats()
$ws.Range("E36").Value = "  -0.08%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.706"
$c.ClearFormats()
$ws.Range("E37").Value = "  -0.43%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01879"
$c.ClearFormats()
$ws.Range("E38").Value = "  +0.75%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.617"
$c.ClearFormats()
$ws.Range("E39").Value = "  +3.04%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.242"
$c.ClearFormats()
$ws.Range("E40").Value = "  -2.60%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "74.60"
$c.ClearFormats()
$ws.Range("E41").Value = "  +1.17%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.944"
$c.ClearFormats()
$ws.Range("E42").Value = "  -0.20%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.8463"
$c.ClearFormats()
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  -0.60%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.9994"
$c.ClearFormats()
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E46").Value = "  -0.86%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "956.00"
$c.ClearFormats()
$ws.Range("E47").Value = "  +2.63%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.174"
$c.ClearFormats()
$ws.Range("E48").Value = "  -0.83%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "9.230"
$c.ClearFormats()
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("E50").Value = "  -0.03%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.05660"
$c.ClearFormats()
$ws.Range("E51").Value = "  +0.43%  "
